$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7216.6665
$ws.Range("H98").Value = 3355
$ws.Range("I98").Value = 3355
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 3355
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -1857
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 3355
$ws.Range("I122").Value = 3355
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10065
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7615
$ws.Range("N122").ClearContents()
$ws.Range("H138").Value = 3287.682
$ws.Range("I138").Value = 871
$ws.Range("J138").Value = 3998.4707
$ws.Range("K138").Value = 2613
$ws.Range("L138").Value = 11995.4121
$ws.Range("M138").Value = 2527
$ws.Range("N138").Value = -22275.4121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 3497.5
$ws.Range("I10").Value = 3000
$ws.Range("K10").Value = 3000
$ws.Range("M10").Value = -2830
$ws.Range("H32").Value = 3385.7778
$ws.Range("I32").Value = 3385.7778
$ws.Range("K32").Value = 3385.7778
$ws.Range("M32").Value = -3098.7778
$ws.Range("H45").Value = 1522.2
$ws.Range("I45").Value = 1522.2
$ws.Range("K45").Value = 1522.2
$ws.Range("M45").Value = -1145.2
$ws.Range("H61").Value = 3354.3635
$ws.Range("I61").Value = 3001.5715
$ws.Range("J61").Value = 3971.75
$ws.Range("K61").Value = 3001.5715
$ws.Range("L61").Value = 3971.75
$ws.Range("M61").Value = -2789.5715
$ws.Range("N61").Value = -4395.75
$ws.Range("H74").Value = 5933.3335
$ws.Range("I74").Value = 5933.3335
$ws.Range("K74").Value = 5933.3335
$ws.Range("M74").Value = -5059.3335
$ws.Range("H77").Value = 5933.3335
$ws.Range("I77").Value = 5933.3335
$ws.Range("K77").Value = 29666.6675
$ws.Range("M77").Value = -25298.6675
$ws.Range("H110").Value = 901.5
$ws.Range("I110").Value = 901.5
$ws.Range("K110").Value = 901.5
$ws.Range("M110").Value = 1143.5
$ws.Range("H122").Value = 6812.4
$ws.Range("I122").Value = 7141.375
$ws.Range("K122").Value = 21424.125
$ws.Range("M122").Value = -18974.125
$ws.Range("H132").Value = 2799.6667
$ws.Range("I132").Value = 2799.6667
$ws.Range("K132").Value = 8399.000100000001
$ws.Range("M132").Value = -5869.000100000001
$ws.Range("H136").Value = 3354.3635
$ws.Range("I136").Value = 3001.5715
$ws.Range("J136").Value = 3971.75
$ws.Range("K136").Value = 9004.7145
$ws.Range("L136").Value = 11915.25
$ws.Range("M136").Value = -6454.7145
$ws.Range("N136").Value = -17015.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 550.75
$ws.Range("J5").Value = 517.5
$ws.Range("L5").Value = 517.5
$ws.Range("N5").Value = -743.5
$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4664
$ws.Range("H134").Value = 6300.4
$ws.Range("I134").Value = 806
$ws.Range("J134").Value = 11794.8
$ws.Range("K134").Value = 2418
$ws.Range("L134").Value = 35384.39999999999
$ws.Range("M134").Value = 117
$ws.Range("N134").Value = -40454.39999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6003.643
$ws.Range("I31").Value = 4926.6665
$ws.Range("K31").Value = 4926.6665
$ws.Range("M31").Value = -4631.6665
$ws.Range("H34").Value = 6003.643
$ws.Range("I34").Value = 4926.6665
$ws.Range("K34").Value = 4926.6665
$ws.Range("M34").Value = -4724.6665
$ws.Range("H51").Value = 22444
$ws.Range("J51").Value = 22444
$ws.Range("L51").Value = 22444
$ws.Range("N51").Value = -23916
$ws.Range("H58").Value = 312.5
$ws.Range("I58").Value = 312.5
$ws.Range("K58").Value = 312.5
$ws.Range("M58").Value = -109.5
$ws.Range("H59").Value = 35000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H61").Value = 22444
$ws.Range("J61").Value = 22444
$ws.Range("L61").Value = 22444
$ws.Range("N61").Value = -23140
$ws.Range("H62").Value = 3510.5
$ws.Range("I62").Value = 3347.3333
$ws.Range("K62").Value = 3347.3333
$ws.Range("M62").Value = -2723.3333
$ws.Range("H65").Value = 3510.5
$ws.Range("I65").Value = 3347.3333
$ws.Range("K65").Value = 16736.6665
$ws.Range("M65").Value = -13616.6665
$ws.Range("H136").Value = 312.5
$ws.Range("I136").Value = 312.5
$ws.Range("K136").Value = 937.5
$ws.Range("M136").Value = 1612.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.909092
$ws.Range("I2").Value = 17.5
$ws.Range("K2").Value = 105
$ws.Range("M2").Value = 8
$ws.Range("H11").Value = 489.875
$ws.Range("I11").Value = 486.83334
$ws.Range("K11").Value = 1460.50002
$ws.Range("M11").Value = -1320.50002
$ws.Range("H17").Value = 262.5
$ws.Range("J17").Value = 262.5
$ws.Range("L17").Value = 787.5
$ws.Range("N17").Value = -1125.5
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H139").Value = 3369.4443
$ws.Range("I139").Value = 2554.1667
$ws.Range("K139").Value = 7662.500100000001
$ws.Range("M139").Value = -2522.500100000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 874.625
$ws.Range("I102").Value = 874.625
$ws.Range("K102").Value = 874.625
$ws.Range("M102").Value = 747.375
$ws.Range("H132").Value = 2676.4707
$ws.Range("I132").Value = 1928.6364
$ws.Range("J132").Value = 4047.5
$ws.Range("K132").Value = 5785.9092
$ws.Range("L132").Value = 12142.5
$ws.Range("M132").Value = -3255.9092
$ws.Range("N132").Value = -17202.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1350
$ws.Range("I32").Value = 1350
$ws.Range("K32").Value = 1350
$ws.Range("M32").Value = -1033
$ws.Range("H46").Value = 2808.5
$ws.Range("I46").Value = 1198.8
$ws.Range("J46").Value = 3958.2856
$ws.Range("K46").Value = 1198.8
$ws.Range("L46").Value = 3958.2856
$ws.Range("M46").Value = -1010.8
$ws.Range("N46").Value = -4334.2856
$ws.Range("H55").Value = 2279.5
$ws.Range("I55").Value = 250.28572
$ws.Range("J55").Value = 3857.7778
$ws.Range("K55").Value = 250.28572
$ws.Range("L55").Value = 3857.7778
$ws.Range("M55").Value = -77.28572
$ws.Range("N55").Value = -4203.7778
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 56000
$ws.Range("J130").Value = 56000
$ws.Range("L130").Value = 56000
$ws.Range("N130").Value = -66040
